$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.004.19"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.641.06"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.54"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5151"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06373"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.290"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.641.54"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5478"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "0.0₅7787"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.55"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "26.027.88"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "199.27"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.464"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.988"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.098"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.907"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.31"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1237"
$ws.Range("E26").Value = "  +8.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.876"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04870"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.308"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.231"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.378"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9183"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.154.04"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5596"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.573"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01575"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.587"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8088"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.79"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").Value = "1.781.67"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.35"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05215"
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09581"
$ws.Range("E51").Value = "  +0.27%  "
